$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily row appended at the bottom of the log (row 76).
# The date column holds a date-formatted string ("2025/10/07"); a leading
# apostrophe forces Excel to store it as literal text instead of auto
# converting it to a date serial, and resetting the style back to "Normal"
# afterwards keeps the cell on the sheet's default (unstyled) format, same
# as every other data row.
$ws.Cells.Item(76, 1).Value = "'2025/10/07"
$ws.Cells.Item(76, 1).Style = "Normal"

$ws.Cells.Item(76, 2).Value = "火"
$ws.Cells.Item(76, 3).Value = 22
$ws.Cells.Item(76, 4).Value = 201
